# daily auto push: 2026-02-10 10:08 UTC
# Inserts a new daily record row at row 778 (date 2026/02/10, time 17, value 151),
# shifting the existing rows 778:819 down to 779:820.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 778 - shifts rows 778..819 down to 779..820
$ws.Rows.Item(778).Insert()

# Populate the newly inserted row with the new daily record.
# Force column A to be treated as text first so the "yyyy/mm/dd" string
# isn't auto-converted into a date serial number, then strip the
# temporary number format back off so the cell keeps the sheet's
# default (unstyled) look, matching the rest of the date column.
$ws.Range("A778").NumberFormat = "@"
$ws.Range("A778").Value = "2026/02/10"
$ws.Range("A778").ClearFormats()

$ws.Range("B778").Value = "火"
$ws.Range("C778").Value = 17
$ws.Range("D778").Value = 151
